$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header labels to reflect "six year" terminology instead of "150%"
$ws.Range("C1").Value = "Male Completers six years"
$ws.Range("E1").Value = "Female Completers six years"
$ws.Range("H1").Value = "six year graduation count"

# Match the cursor/selection position left behind in the source file
$ws.Range("F8").Select()
